$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds free-form price text (e.g. '26.041.68', '18.61').
# Force the whole data range to Text first so the COM layer doesn't
# coerce plain-looking numerics (like '18.61') into floating point
# numbers when we assign them below.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.041.68'
$ws.Range("E2").Value = '  +0.06%  '
$ws.Range("D3").Value = '1.630.64'
$ws.Range("E3").Value = '  -0.85%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '214.34'
$ws.Range("E5").Value = '  -0.53%  '
$ws.Range("E8").Value = '  -1.78%  '
$ws.Range("E9").Value = '  -2.94%  '
$ws.Range("D10").Value = '18.61'
$ws.Range("E10").Value = '  -5.03%  '
$ws.Range("E11").Value = '  -0.94%  '
$ws.Range("D12").Value = '1.857.17'
$ws.Range("E12").Value = '  -0.86%  '
$ws.Range("D13").Value = '1.639.71'
$ws.Range("E13").Value = '  -0.10%  '
$ws.Range("E14").Value = '  -1.73%  '
$ws.Range("D15").Value = '0.528'
$ws.Range("D16").Value = '26.041.16'
$ws.Range("E16").Value = '  -0.12%  '
$ws.Range("D17").Value = '0.0₃0743'
$ws.Range("E17").Value = '  -2.31%  '
$ws.Range("D18").Value = '61.63'
$ws.Range("E18").Value = '  -2.88%  '
$ws.Range("E19").Value = '  -0.01%  '
$ws.Range("D20").Value = '193.06'
$ws.Range("E20").Value = '  -0.65%  '
$ws.Range("E21").Value = '  -2.15%  '
$ws.Range("E22").Value = '  -3.52%  '
$ws.Range("E24").Value = '  +1.97%  '
$ws.Range("D25").Value = '144.31'
$ws.Range("E25").Value = '  +0.23%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("E27").Value = '  -3.74%  '
$ws.Range("E28").Value = '  -2.14%  '
$ws.Range("D29").Value = '15.28'
$ws.Range("E29").Value = '  -1.39%  '
$ws.Range("E30").Value = '  -0.72%  '
$ws.Range("E31").Value = '  -2.52%  '
$ws.Range("D32").Value = '3.13'
$ws.Range("E32").Value = '  -3.95%  '
$ws.Range("E33").Value = '  -4.71%  '
$ws.Range("E34").Value = '  -3.01%  '
$ws.Range("E35").Value = '  -1.94%  '
$ws.Range("D36").Value = '1.128.73'
$ws.Range("E36").Value = '  -0.19%  '
$ws.Range("E37").Value = '  -5.56%  '
$ws.Range("E38").Value = '  -1.19%  '
$ws.Range("E39").Value = '  -3.25%  '
$ws.Range("E40").Value = '  -2.23%  '
$ws.Range("D41").Value = '98.31'
$ws.Range("E41").Value = '  -0.63%  '
$ws.Range("D42").Value = '1.766.67'
$ws.Range("E42").Value = '  -0.87%  '
$ws.Range("D43").Value = '0.759'
$ws.Range("E43").Value = '  -4.67%  '
$ws.Range("D44").Value = '5.13'
$ws.Range("E44").Value = '  -5.71%  '
$ws.Range("E45").Value = '  -1.86%  '
$ws.Range("D46").Value = '54.59'
$ws.Range("E46").Value = '  -3.37%  '
$ws.Range("D47").Value = '0.0523'
$ws.Range("E47").Value = '  +0.27%  '
$ws.Range("E48").Value = '  -0.18%  '
$ws.Range("E49").Value = '  -0.27%  '
$ws.Range("D50").Value = '7.50'
$ws.Range("E50").Value = '  -3.70%  '
$ws.Range("E51").Value = '  +0.02%  '

# Restore original (unstyled) formatting now that the text values are
# committed as strings.
$ws.Range("D2:D51").ClearFormats()
